$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a long block of pre-styled but empty rows (29-125) reserved
# for future diary entries. This edit fills in three more entries (rows
# 29-31) describing work done up through 2020-03-14. Copy the formatting
# from the last filled-in data row (28) down onto the three rows we are
# about to populate, matching the look of the other entries.
$ws.Range("A28:G28").Copy()
$ws.Range("A29:G31").PasteSpecial(-4122)

# Row 29 - 2020-03-05
$ws.Cells.Item(29, 1).Value = 43895
$ws.Cells.Item(29, 2).Value = "5:00PM - 7:50 PM"
$ws.Cells.Item(29, 3).Value = "N/A"
$ws.Cells.Item(29, 4).Value = "Understand how test cases can help us understand code"
$ws.Cells.Item(29, 5).Value = "Learned new key expert practices and how test cases can be useful to actually understand how the code is supposed to behave"
$ws.Cells.Item(29, 6).Value = "Test cases can show us through a input, output format what we should expect from certain parts of the code. If we map that information and associate it with domain knowledge we can learn valuable insight."
$ws.Cells.Item(29, 7).Value = "Feeling good overall"

# Row 30 - 2020-03-12
$ws.Cells.Item(30, 1).Value = 43902
$ws.Cells.Item(30, 2).Value = "5:00PM - 7:50 PM"
$ws.Cells.Item(30, 3).Value = "N/A"
$ws.Cells.Item(30, 4).Value = "Learn advanced topics"
$ws.Cells.Item(30, 5).Value = "Learned new key expert practices, how the history of the project can be relevant, and how visualizations can offer good insight"
$ws.Cells.Item(30, 6).Value = "Even after the program is over, I must keep studying to stay up to date with the current trends. This way I am almost always ready for what is coming"
$ws.Cells.Item(30, 7).Value = "Feeling good, but a bit tired"

# Row 31 - 2020-03-14
$ws.Cells.Item(31, 1).Value = 43904
$ws.Cells.Item(31, 2).Value = "2:00PM - 7:00 PM"
$ws.Cells.Item(31, 3).Value = "Chris Zhang, Nicolas Grantham, and Hyun Jay Yang"
$ws.Cells.Item(31, 4).Value = "Finish the last assignment for the class"
$ws.Cells.Item(31, 5).Value = "Contributed with our second issue, studying and detailed three test cases that we found interesting, and created new test cases for our project"
$ws.Cells.Item(31, 6).Value = "Having good knowledge of the data-flow, control-flow, and architecture of our project made our contribution so much easier, since we knew what had to be changed, making the code understanding part of the assignment easier."
$ws.Cells.Item(31, 7).Value = "Feeling tired, but glad to have contributed"

# Move the active selection to reflect where the author's cursor ended up.
$ws.Range("G30").Select()
